$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 78: fx / sftgif variable (traced 126-table variable combination)
$ws.Range("A78").Value = "fx"
$ws.Range("B78").Value = "sftgif"
$ws.Range("C78").NumberFormat = "@"
$ws.Range("C78").Value = "1"
$ws.Range("C78").NumberFormat = "General"
$ws.Range("D78").Value = "longitude latitude time typeli"
$ws.Range("E78").Value = "Land Ice Area Percentage"
$ws.Range("F78").Value = "%"
$ws.Range("G78").Formula = '=HYPERLINK("http://clipc-services.ceda.ac.uk/dreq/u/a1d2e309c6f25017442ad6c79c4f9eca.html","web")'
$ws.Range("H78").Value = "To be implemented:  grib 126.32  part of MFPPHY   This is the land ice mask and will be an extra variable in IFS (thomas: via PEXTRA?)"
$ws.Range("I78").Value = "Shuting"
$ws.Range("J78").Value = "Fraction of grid cell covered by land ice (ice sheet, ice shelf, ice cap, glacier)"
$ws.Range("K78").Value = "CMIP,ISMIP6"

# Row 80: Emon / vtendogw
$ws.Range("A80").Value = "Emon"
$ws.Range("B80").Value = "vtendogw"
$ws.Range("D80").Value = "longitude latitude plev19 time"
$ws.Range("E80").Value = "Northward Acceleration Due to Orographic Gravity Wave Drag"
$ws.Range("F80").Value = "m s-2"
$ws.Range("G80").Value = 0
$ws.Range("J80").Value = "Tendency of the northward wind by parameterized orographic gravity waves.  (Note that CF name tables only have a general northward tendency for all gravity waves, and we need it separated by type.)"

# Row 81: Emon / vtendnogw
$ws.Range("A81").Value = "Emon"
$ws.Range("B81").Value = "vtendnogw"
$ws.Range("D81").Value = "longitude latitude plev19 time"
$ws.Range("E81").Value = "Northward Acceleration Due to Non-Orographic Gravity Wave Drag"
$ws.Range("F81").Value = "m s-2"
$ws.Range("G81").Value = 0
$ws.Range("J81").Value = "Tendency of the northward wind by parameterized nonorographic gravity waves.  (Note that CF name tables only have a general northward tendency for all gravity waves, and we need it separated by type.)"

# Row 82: EmonZ / vtendnogw
$ws.Range("A82").Value = "EmonZ"
$ws.Range("B82").Value = "vtendnogw"
$ws.Range("D82").Value = "longitude latitude plev19 time"
$ws.Range("E82").Value = "Northward Acceleration Due to Non-Orographic Gravity Wave Drag"
$ws.Range("F82").Value = "m s-2"
$ws.Range("G82").Value = 0
$ws.Range("J82").Value = "Tendency of the northward wind by parameterized nonorographic gravity waves.  (Note that CF name tables only have a general northward tendency for all gravity waves, and we need it separated by type.)"

# Row 83: EmonZ / tntogw
$ws.Range("A83").Value = "EmonZ"
$ws.Range("B83").Value = "tntogw"
$ws.Range("D83").Value = "latitude plev39 time"
$ws.Range("E83").Value = "Temperature Tendency Due to Orographic Gravity Wave Dissipation"
$ws.Range("F83").Value = "K s-1"
$ws.Range("G83").Value = 0
$ws.Range("J83").Value = "Temperature tendency due to dissipation of parameterized orographic gravity waves."

# Row 84: EmonZ / tntnogw
$ws.Range("A84").Value = "EmonZ"
$ws.Range("B84").Value = "tntnogw"
$ws.Range("D84").Value = "latitude plev39 time"
$ws.Range("E84").Value = "Temperature Tendency Due to Non-Orographic Gravity Wave Dissipation"
$ws.Range("F84").Value = "K s-1"
$ws.Range("G84").Value = 0
$ws.Range("J84").Value = "Temperature tendency due to dissipation of parameterized nonorographic gravity waves."

# Update the view: scroll so row 59 is the top-left visible row, and select C80:C84
$excel.ActiveWindow.ScrollRow = 59
$ws.Range("C80:C84").Select()
